$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row (row 1) - reorder block labels
$ws.Range("B1").Value = "kitchens_2"
$ws.Range("D1").Value = "bedrooms_1"
$ws.Range("E1").Value = "bedrooms_2"
$ws.Range("F1").Value = "living_rooms_2"

# Update data rows to match the new column order
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0

$ws.Range("B3").Value = 1
$ws.Range("D3").Value = 0

$ws.Range("E4").Value = 0
$ws.Range("F4").Value = 1

$ws.Range("B5").Value = 0
$ws.Range("D5").Value = 1
